$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 ---
$ws.Range("A2").Value2 = 131223149
$ws.Range("B2").Value2 = 57884
$ws.Range("E2").Value2 = 100109
$ws.Range("F2").Value2 = "Tretåig hackspett"
$ws.Range("G2").Value2 = "Picoides tridactylus"
$ws.Range("H2").Value2 = "(Linnaeus, 1758)"
$ws.Range("J2").ClearContents()
$ws.Range("L2").Style = "Normal"
$ws.Range("M2").Value2 = "färska spår"
$ws.Range("Q2").Value2 = 772974
$ws.Range("R2").Value2 = 7122563
$ws.Range("AC2").Value2 = "färska ringhack på gran"
$ws.Range("AF2").ClearContents()

# --- Row 3 ---
$ws.Range("A3").Value2 = 131223060
$ws.Range("Q3").Value2 = 772981
$ws.Range("R3").Value2 = 7122639

# --- Row 4 ---
$ws.Range("A4").Value2 = 131223489
$ws.Range("B4").Value2 = 79245
$ws.Range("E4").Value2 = 6425
$ws.Range("F4").Value2 = "Garnlav"
$ws.Range("G4").Value2 = "Alectoria sarmentosa"
$ws.Range("H4").Value2 = "(Ach.) Ach."
$ws.Range("J4").Value2 = "bålar"
$ws.Range("L4").ClearContents()
$ws.Range("M4").ClearContents()
$ws.Range("Q4").Value2 = 773011
$ws.Range("R4").Value2 = 7122664
$ws.Range("AC4").ClearContents()
$ws.Range("AF4").Style = "Normal"

# --- Row 16 ---
$ws.Range("B16").Value2 = 91810
